try {
  [System.Environment]::SetEnvironmentVariable("IRON_VBA_MAC_PLATFORM", "1", "Process")
  Write-Host "set done"
} catch {
  Write-Host "ERR: $_"
}
Write-Host "Value: $([System.Environment]::GetEnvironmentVariable('IRON_VBA_MAC_PLATFORM'))"
